$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.406.61'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').Value = '3.220.08'
$ws.Range('E3').Value = '  -1.55%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '577.23'
$ws.Range('E5').Value = '  -1.73%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '182.39'
$ws.Range('E6').Value = '  -1.18%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('D9').Value = '3.218.52'
$ws.Range('E9').Value = '  -1.49%  '
$ws.Range('E10').Value = '  -3.41%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '6.57'
$ws.Range('E11').Value = '  -2.10%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.411'
$ws.Range('E12').Value = '  -1.49%  '
$ws.Range('D13').Value = '3.778.28'
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '27.64'
$ws.Range('E15').Value = '  -3.66%  '
$ws.Range('D16').Value = '67.456.06'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('E17').Value = '  -2.46%  '
$ws.Range('D18').Value = '3.237.72'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('E19').Value = '  -2.29%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.38'
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '393.88'
$ws.Range('E21').Value = '  +2.75%  '
$ws.Range('E22').Value = '  -2.42%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.33%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '70.81'
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.513'
$ws.Range('E25').Value = '  -0.61%  '
$ws.Range('E26').Value = '  -3.50%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.184'
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.53'
$ws.Range('E28').Value = '  -3.57%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  -2.50%  '
$ws.Range('E31').Value = '  -3.76%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '22.57'
$ws.Range('E32').Value = '  -1.63%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.94'
$ws.Range('E33').Value = '  -4.38%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  -2.82%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '160.72'
$ws.Range('E36').Value = '  -1.29%  '
$ws.Range('E37').Value = '  -5.53%  '
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('E39').Value = '  -3.99%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '26.12'
$ws.Range('E40').Value = '  -2.42%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '4.54'
$ws.Range('E41').Value = '  -1.64%  '
$ws.Range('E42').Value = '  -4.32%  '
$ws.Range('E43').Value = '  -5.88%  '
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '40.48'
$ws.Range('E45').Value = '  -2.46%  '
$ws.Range('D46').Value = '2.587.02'
$ws.Range('E46').Value = '  -2.46%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '24.43'
$ws.Range('E47').Value = '  -4.09%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '332.08'
$ws.Range('E48').Value = '  -4.89%  '
$ws.Range('E49').Value = '  -3.04%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.26'
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('E51').Value = '  -2.03%  '
